$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (SVM)
$ws.Range("B2").Value = 0.733
$ws.Range("C2").Value = 0.001
$ws.Range("D2").Value = 0.787
$ws.Range("F2").Value = 0.0005446055589503767

# Row 3 (LR)
$ws.Range("B3").Value = 0.721
$ws.Range("D3").Value = 0.784

# Row 4 (LDA)
$ws.Range("B4").Value = 0.664
$ws.Range("D4").Value = 0.693

# Row 5 (RF)
$ws.Range("B5").Value = 0.617
$ws.Range("C5").Value = 0.009
$ws.Range("D5").Value = 0.647
$ws.Range("F5").Value = 0.008722348819355559

# Row 6 (AB)
$ws.Range("B6").Value = 0.649
$ws.Range("D6").Value = 0.687

# Row 7 (KNN)
$ws.Range("B7").Value = 0.631
$ws.Range("D7").Value = 0.654

# Row 8 (GNB)
$ws.Range("B8").Value = 0.644
$ws.Range("D8").Value = 0.681
